# Danh sách khách hàng tại SÓC TRĂNG - thêm 5 khách hàng mới vào đầu danh sách.
# fix lỗi trong report cơ sở. Thêm cột ghi chú trong báo cáo về chi tiêu

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows right after the header row (row 1), pushing all
# existing customer rows down by 5.
$ws.Rows("2:6").Insert()

# New customer records (columns: A Tiền tố, B Mã khách hàng, C Họ và tên,
# D Cơ sở, E CCCD, F SĐT, G Facebook, H Địa chỉ, I Tích lũy, J Dư nợ).
$newCustomers = @(
    @{ Row = 2; B = 411; C = "thạch thị siêu";      F = "0833525730";  I = 50000000; J = 26000000 },
    @{ Row = 3; B = 405; C = "tạ duy hoàng ";        F = $null;         I = 6000000;  J = 0 },
    @{ Row = 4; B = 403; C = "nguyễn thị lệ trang";  F = "0786070121";  I = 4000000;  J = 0 },
    @{ Row = 5; B = 402; C = "nguyễn thị mỹ trinh";  F = $null;         I = 7000000;  J = 0 },
    @{ Row = 6; B = 401; C = "nguyễn thị mỹ chăm";   F = "09399259920"; I = 6000000;  J = 0 }
)

foreach ($cust in $newCustomers) {
    $r = $cust.Row

    $ws.Cells.Item($r, 1).Value = "KH"
    $ws.Cells.Item($r, 2).Value = $cust.B
    $ws.Cells.Item($r, 3).Value = $cust.C
    $ws.Cells.Item($r, 4).Value = "SÓC TRĂNG"

    if ($cust.F) {
        # Store phone numbers as text so leading zeros are preserved.
        $ws.Cells.Item($r, 6).NumberFormat = "@"
        $ws.Cells.Item($r, 6).Value = $cust.F
    }

    $ws.Cells.Item($r, 9).Value = $cust.I
    $ws.Cells.Item($r, 10).Value = $cust.J
}
